$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '28.484.38'
$c.Style = 'Normal'

$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  +0.43%  '
$c.Style = 'Normal'

$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '1.868.92'
$c.Style = 'Normal'

$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +0.14%  '
$c.Style = 'Normal'

$c = $ws.Range("D4")
$c.NumberFormat = '@'
$c.Value = '1.008'
$c.Style = 'Normal'

$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  -1.26%  '
$c.Style = 'Normal'

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '315.49'
$c.Style = 'Normal'

$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.Style = 'Normal'

$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -0.98%  '
$c.Style = 'Normal'

$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.5072'
$c.Style = 'Normal'

$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  -0.81%  '
$c.Style = 'Normal'

$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.3898'
$c.Style = 'Normal'

$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  -0.55%  '
$c.Style = 'Normal'

$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.08350'
$c.Style = 'Normal'

$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +0.88%  '
$c.Style = 'Normal'

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '1.105'
$c.Style = 'Normal'

$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  -0.48%  '
$c.Style = 'Normal'

$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '41.78'
$c.Style = 'Normal'

$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  -0.11%  '
$c.Style = 'Normal'

$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '6.217'
$c.Style = 'Normal'

$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  -0.08%  '
$c.Style = 'Normal'

$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '1.874.51'
$c.Style = 'Normal'

$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.Style = 'Normal'

$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  +0.25%  '
$c.Style = 'Normal'

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '7.275'
$c.Style = 'Normal'

$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  +1.05%  '
$c.Style = 'Normal'

$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '1.007'
$c.Style = 'Normal'

$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -1.25%  '
$c.Style = 'Normal'

$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  -0.15%  '
$c.Style = 'Normal'

$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '91.05'
$c.Style = 'Normal'

$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c.Style = 'Normal'

$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '0.06730'
$c.Style = 'Normal'

$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.Style = 'Normal'

$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '17.72'
$c.Style = 'Normal'

$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  +0.54%  '
$c.Style = 'Normal'

$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  -1.01%  '
$c.Style = 'Normal'

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '5.915'
$c.Style = 'Normal'

$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  -0.61%  '
$c.Style = 'Normal'

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '28.509.27'
$c.Style = 'Normal'

$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  +0.46%  '
$c.Style = 'Normal'

$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '11.08'
$c.Style = 'Normal'

$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  -0.20%  '
$c.Style = 'Normal'

$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '2.208'
$c.Style = 'Normal'

$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  -1.80%  '
$c.Style = 'Normal'

$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '2.085.03'
$c.Style = 'Normal'

$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.Style = 'Normal'

$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '158.50'
$c.Style = 'Normal'

$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  -1.50%  '
$c.Style = 'Normal'

$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '20.59'
$c.Style = 'Normal'

$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  -0.52%  '
$c.Style = 'Normal'

$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '2.425'
$c.Style = 'Normal'

$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +2.65%  '
$c.Style = 'Normal'

$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '126.34'
$c.Style = 'Normal'

$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -0.44%  '
$c.Style = 'Normal'

$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  -0.94%  '
$c.Style = 'Normal'

$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  +1.07%  '
$c.Style = 'Normal'

$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '5.734'
$c.Style = 'Normal'

$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  -0.93%  '
$c.Style = 'Normal'

$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '3.630'
$c.Style = 'Normal'

$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  -0.12%  '
$c.Style = 'Normal'

$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '0.02455'
$c.Style = 'Normal'

$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  +1.15%  '
$c.Style = 'Normal'

$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '0.06578'
$c.Style = 'Normal'

$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +1.47%  '
$c.Style = 'Normal'

$c = $ws.Range("B37")
$c.NumberFormat = '@'
$c.Value = 'Algorand'
$c.Style = 'Normal'

$c = $ws.Range("C37")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.Style = 'Normal'

$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '0.2165'
$c.Style = 'Normal'

$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -0.45%  '
$c.Style = 'Normal'

$c = $ws.Range("B38")
$c.NumberFormat = '@'
$c.Value = 'FraxShare'
$c.Style = 'Normal'

$c = $ws.Range("C38")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c.Style = 'Normal'

$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '8.892'
$c.Style = 'Normal'

$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  -2.52%  '
$c.Style = 'Normal'

$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '5.030'
$c.Style = 'Normal'

$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  +1.16%  '
$c.Style = 'Normal'

$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '1.179'
$c.Style = 'Normal'

$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -0.23%  '
$c.Style = 'Normal'

$c = $ws.Range("B41")
$c.NumberFormat = '@'
$c.Value = 'TrustWalletToken'
$c.Style = 'Normal'

$c = $ws.Range("C41")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.Style = 'Normal'

$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '1.233'
$c.Style = 'Normal'

$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  -1.07%  '
$c.Style = 'Normal'

$c = $ws.Range("B42")
$c.NumberFormat = '@'
$c.Value = 'TheSandbox'
$c.Style = 'Normal'

$c = $ws.Range("C42")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c.Style = 'Normal'

$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.6365'
$c.Style = 'Normal'

$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  -0.67%  '
$c.Style = 'Normal'

$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '11.09'
$c.Style = 'Normal'

$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  -0.30%  '
$c.Style = 'Normal'

$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -0.97%  '
$c.Style = 'Normal'

$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.6007'
$c.Style = 'Normal'

$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.Style = 'Normal'

$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '12.99'
$c.Style = 'Normal'

$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  +0.57%  '
$c.Style = 'Normal'

$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.Style = 'Normal'

$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '2.003'
$c.Style = 'Normal'

$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  +1.03%  '
$c.Style = 'Normal'

$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '1.214'
$c.Style = 'Normal'

$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  +0.80%  '
$c.Style = 'Normal'

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '122.36'
$c.Style = 'Normal'

$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  +0.67%  '
$c.Style = 'Normal'

$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '0.06805'
$c.Style = 'Normal'

$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c.Style = 'Normal'
